# Update Structurizer, Utilizer, main, prompts
#
# Inserts 8 new boolean "value/type" column pairs (涉及共犯, 涉及外國人,
# 和解, 被害人考量) ahead of the existing 媒體影響/量刑爭議 pair columns,
# fills in their header labels + row data, and tweaks a handful of other
# cells (L2, V3, Z2, Z3) to match the updated extraction results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 8 new columns starting at AA -------------------
# This shifts the old AA:AD (媒體影響_value/_type, 量刑爭議_value/_type)
# over to AI:AL, preserving their header styling and row values.
$ws.Range("AA1:AH1").EntireColumn.Insert()

# --- 2. New header row labels for the freshly inserted columns -----------
$ws.Range("AA1").Value = "涉及共犯_value"
$ws.Range("AB1").Value = "涉及共犯_type"
$ws.Range("AC1").Value = "涉及外國人_value"
$ws.Range("AD1").Value = "涉及外國人_type"
$ws.Range("AE1").Value = "和解_value"
$ws.Range("AF1").Value = "和解_type"
$ws.Range("AG1").Value = "被害人考量_value"
$ws.Range("AH1").Value = "被害人考量_type"

# --- 3. Row 2 data ---------------------------------------------------------
$ws.Range("L2").Value = $false

$ws.Range("AA2").Value = $false
$ws.Range("AB2").Value = "bool"
$ws.Range("AC2").Value = $false
$ws.Range("AD2").Value = "bool"
$ws.Range("AE2").Value = $false
$ws.Range("AF2").Value = "bool"
$ws.Range("AG2").Value = $true
$ws.Range("AH2").Value = "bool"
$ws.Range("AI2").Value = $true
$ws.Range("AJ2").Value = "bool"

$ws.Range("Z2").Value = "因為被告已對被訴事實作有罪陳述，且案件情節涉及媒體影響及被害人家屬的意見，法院認為行國民參與審判可能對被害人家屬造成二度傷害，且不適合彰顯國民參與審判的價值。"

# --- 4. Row 3 data ---------------------------------------------------------
$ws.Range("V3").Value = $true

$ws.Range("AA3").Value = $false
$ws.Range("AB3").Value = "bool"
$ws.Range("AC3").Value = $false
$ws.Range("AD3").Value = "bool"
$ws.Range("AE3").Value = $true
$ws.Range("AF3").Value = "bool"
$ws.Range("AG3").Value = $true
$ws.Range("AH3").Value = "bool"
$ws.Range("AK3").Value = $false

$ws.Range("Z3").Value = "因為被告已對被訴事實為有罪陳述，且檢辯雙方對於量刑無重大爭議，並且被害人家屬已表達同意不進行國民參與審判程序，符合國民法官法第6條第1項第4款的規定。"
